$d = $word.ActiveDocument

# Locate the existing "DB AUTO USATE" run (the document title) and
# collapse the found range to its start so we can insert a brand-new run
# right before it, rather than merging text into the existing run.
$find = $d.Content.Find
$find.Text = "DB AUTO USATE"
$find.Execute() | Out-Null
$r = $find.Parent
$r.Collapse(1)

# Insert a new run containing "202" with formatting identical to the
# "DB AUTO USATE" run (Helvetica Neue, bold, sz 34 / szCs 34). Using
# InsertXML (instead of InsertBefore / Range.Text) preserves this as a
# distinct <w:r> element rather than being coalesced into the
# neighbouring run's text.
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData>' + `
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:body><w:p><w:r><w:rPr>' + `
    '<w:rFonts w:ascii="Helvetica Neue" w:hAnsi="Helvetica Neue" w:cs="Helvetica Neue"/>' + `
    '<w:b/><w:bCs/><w:sz w:val="34"/><w:szCs w:val="34"/>' + `
    '</w:rPr><w:t>202</w:t></w:r></w:p></w:body>' + `
    '</w:document>' + `
    '</pkg:xmlData></pkg:part></pkg:package>'

$r.InsertXML($xml)
